$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial number (45180 -> 2023-09-11)
# that is bumped by one day (45181 -> 2023-09-12) for every data row.
$ws.Range("C2:C480").Value = 45181
